# Update the "scraped_at" timestamps (column K) on the "snapshot" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("snapshot")

$timestamps = @{
    2  = "2025-12-20T13:51:17.095633+00:00"
    3  = "2025-12-20T13:51:19.137637+00:00"
    4  = "2025-12-20T13:51:19.137655+00:00"
    5  = "2025-12-20T13:51:19.137663+00:00"
    6  = "2025-12-20T13:51:21.159758+00:00"
    7  = "2025-12-20T13:51:23.440802+00:00"
    8  = "2025-12-20T13:51:25.785088+00:00"
    9  = "2025-12-20T13:51:25.785118+00:00"
    10 = "2025-12-20T13:51:28.171871+00:00"
    11 = "2025-12-20T13:51:33.329149+00:00"
    12 = "2025-12-20T13:51:35.635932+00:00"
    13 = "2025-12-20T13:51:38.547075+00:00"
    14 = "2025-12-20T13:51:44.153046+00:00"
    15 = "2025-12-20T13:51:44.153077+00:00"
    16 = "2025-12-20T13:51:44.153096+00:00"
    17 = "2025-12-20T13:51:44.153113+00:00"
    18 = "2025-12-20T13:51:46.491821+00:00"
    19 = "2025-12-20T13:51:46.491855+00:00"
    20 = "2025-12-20T13:51:46.491874+00:00"
    21 = "2025-12-20T13:51:48.940353+00:00"
    22 = "2025-12-20T13:51:48.940383+00:00"
    23 = "2025-12-20T13:51:48.940403+00:00"
    24 = "2025-12-20T13:51:51.791177+00:00"
    25 = "2025-12-20T13:51:51.791209+00:00"
    26 = "2025-12-20T13:51:54.103131+00:00"
    27 = "2025-12-20T13:51:54.103164+00:00"
    28 = "2025-12-20T13:51:54.103185+00:00"
    29 = "2025-12-20T13:51:56.393785+00:00"
    30 = "2025-12-20T13:51:59.233008+00:00"
    31 = "2025-12-20T13:51:59.233035+00:00"
    32 = "2025-12-20T13:52:03.701279+00:00"
    33 = "2025-12-20T13:52:03.701309+00:00"
    34 = "2025-12-20T13:52:06.196679+00:00"
    35 = "2025-12-20T13:52:06.196709+00:00"
}

foreach ($row in $timestamps.Keys) {
    $ws.Cells.Item($row, 11).Value = $timestamps[$row]
}
